$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($table, $row, $col, $text) {
    $table.Cell($row, $col).Range.Text = $text
}

Set-Cell $t 1 1 "68÷3=22, 2"
Set-Cell $t 1 3 "65÷6=10, 5"
Set-Cell $t 1 4 "89÷7=12, 5"
Set-Cell $t 1 5 "31÷9=3, 4"

Set-Cell $t 5 1 "19÷3=6, 1"
Set-Cell $t 5 2 "21÷3=7, 0"
Set-Cell $t 5 3 "11÷2=5, 1"
Set-Cell $t 5 4 "45÷6=7, 3"
Set-Cell $t 5 5 "95÷2=47, 1"

Set-Cell $t 9 1 "78÷9=8, 6"
Set-Cell $t 9 2 "79÷4=19, 3"
Set-Cell $t 9 3 "19÷9=2, 1"
Set-Cell $t 9 4 "76÷2=38, 0"
Set-Cell $t 9 5 "50÷9=5, 5"

Set-Cell $t 13 1 "62÷5=12, 2"
Set-Cell $t 13 2 "97÷8=12, 1"
Set-Cell $t 13 3 "43÷3=14, 1"
Set-Cell $t 13 4 "24÷9=2, 6"
Set-Cell $t 13 5 "69÷3=23, 0"

Set-Cell $t 17 1 "24÷9=2, 6"
Set-Cell $t 17 2 "57÷5=11, 2"
Set-Cell $t 17 3 "13÷4=3, 1"
Set-Cell $t 17 4 "86÷6=14, 2"
Set-Cell $t 17 5 "55÷2=27, 1"

Write-Host "Table cells updated"
